$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.929.42'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '1.666.71'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.63'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.534'
$ws.Range('E6').Value = '  +5.03%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +0.84%  '
$ws.Range('E9').Value = '  -0.35%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.26'
$ws.Range('E10').Value = '  +2.55%  '
$ws.Range('E11').Value = '  +3.87%  '
$ws.Range('D12').Value = '1.900.34'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '1.670.94'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('E14').Value = '  -0.10%  '
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.18'
$ws.Range('E16').Value = '  +1.39%  '
$ws.Range('D17').Value = '26.934.82'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '236.98'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '8.04'
$ws.Range('E19').Value = '  +1.37%  '
$ws.Range('D20').Value = '0.0₃0733'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.35'
$ws.Range('E22').Value = '  -2.11%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.21'
$ws.Range('E23').Value = '  -1.79%  '
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.98'
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.115'
$ws.Range('E27').Value = '  +1.22%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.92'
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('D33').Value = '1.457.95'
$ws.Range('E33').Value = '  -4.74%  '
$ws.Range('E34').Value = '  +2.50%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.63'
$ws.Range('E35').Value = '  +1.12%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.581'
$ws.Range('E37').Value = '  +0.28%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.906'
$ws.Range('E38').Value = '  +1.25%  '
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.72'
$ws.Range('E40').Value = '  -3.91%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.28'
$ws.Range('E42').Value = '  +0.38%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '65.87'
$ws.Range('E43').Value = '  -1.29%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.973'
$ws.Range('E44').Value = '  +6.10%  '
$ws.Range('D45').Value = '1.809.34'
$ws.Range('E45').Value = '  +0.60%  '
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.53'
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.53'
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0105'
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range('E50').Value = '  +4.42%  '
$ws.Range('E51').Value = '  +0.17%  '
